$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: fill in hours worked and activity description
$ws.Range("E16").Value = 0.7
$ws.Range("F16").Value = "Added missing model classes into backend application. Fixed infinite looping object problem."

# Update selection to match the author's final cursor position (F17)
$ws.Range("F17").Select() | Out-Null
